$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Horas insumidas")

# New timesheet entry: Nico worked 1 hour on "Desarrollo Metricas Agentes" for
# user story S-01004 on 2010-10-16 (serial date 40467).
$ws.Cells.Item(72, 2).Value = 40467
$ws.Cells.Item(72, 3).Value = "Nico"
$ws.Cells.Item(72, 4).Value = "Desarrollo Metricas Agentes"
$ws.Cells.Item(72, 5).Value = "S-01004"
$ws.Cells.Item(72, 6).Value = 1

$ws.Range("F73").Select()

$wb.Save()
